# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (Office colours)   [notes master]
#   ppt/theme/theme2.xml -> "Integral"     (Red Violet colours) [slide master / active design]
#
# The target revision swaps those two palettes: the design that is actually
# applied to the slide master (theme2.xml) goes from "Integral"/Red Violet
# back to the stock "Office Theme"/Office palette. Apply that by rewriting
# every theme colour slot (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# through the ThemeColorScheme COM surface, which is the live, edit-capable
# handle onto the slide master's theme part.

$p = $ppt.ActivePresentation

$officeColors = @(6, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
$officeColors[0] = 0

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}
